# Docs: Fix a typo in the product backlog.
# Adds sequential "Priority" numbers (column A) to the remaining backlog
# rows (15-18) to continue the numbering started in rows 9-14, and moves
# the active selection as it was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14's A cell already carries the numbering style used throughout
# column A (centered text, bordered). Reuse its formatting for the new
# cells in A15:A18, then fill in the next numbers in the sequence.
$ws.Range("A14").Copy()

$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 7

$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 8

$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 9

$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = 10

# Restore the selection to where the author left it (B21).
$ws.Range("B21").Select()
